$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D28").Value = "강화학습 실습 1편 : OpenAI GYM-Atari 환경 셋팅 (Window 10, anaconda)"
$ws.Range("E28").Value = "https://ropiens.tistory.com/153"

$ws.Range("D41").Value = "Stored Procedure에 대한 단상"
$ws.Range("E41").Value = "http://cloudinsight.net/data/stored-procedure%ec%97%90-%eb%8c%80%ed%95%9c-%eb%8b%a8%ec%83%81/"
